$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the species placeholder columns to the actual species names
# (order matters for shared-string table ordering: L7 first, then E2/E7)
$ws.Range("L7").Value = "Abyssorchomene charcoti"
$ws.Range("E2").Value = "Charcotia obesa"
$ws.Range("E7").Value = "Charcotia obesa"

# Match formatting of E7 to the other renamed header cells (E2/L7)
$ws.Range("E7").Font.Bold = $true
$ws.Range("E7").Borders.LineStyle = 1

# Resize columns D and E to fit new content
$ws.Columns("D").ColumnWidth = 15.25
$ws.Columns("E").ColumnWidth = 13.25

# Update the active cell selection
$ws.Range("F9").Select() | Out-Null
